$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the bounding box coordinates for the Erlangen row (row 40)
$ws.Range("D40").Value = 10.9153629
$ws.Range("E40").Value = 49.5327088
$ws.Range("F40").Value = 11.0536043
$ws.Range("G40").Value = 49.6455844
